# Autogenerated on Sun Feb 01 2015 22:24:41 GMT-0500 (Eastern Standard Time)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet from "Data" to "Summary" -----------------------------
$ws.Name = "Summary"

# --- Re-assert formatting on the pre-existing header cells (defensive; --
# --- keeps the big "name" title and bold section title intact). --------
$ws.Range("A1").Font.Size = 18
$ws.Range("A3").Font.Bold = $true

# --- Clear the old "Micro/SMEs/MSMEs" table + source line (rows 5-7); -
# --- it is being re-laid-out further down the sheet. -------------------
$ws.Range("A5:D7").Clear()

# --- New "Source Type" sub-heading (bold + underlined) -----------------
$ws.Range("A9").Value = "Source Type: SME Associations (Most Widely Used)"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# --- Column headers for the MSME participation table --------------------
$ws.Range("B11").Value = "Micro"
$ws.Range("B11").Font.Bold = $true

$ws.Range("C11").Value = "SMEs"
$ws.Range("C11").Font.Bold = $true

$ws.Range("D11").Value = "MSMEs"
$ws.Range("D11").Font.Bold = $true

# --- Employment row (new) ------------------------------------------------
$ws.Range("A12").Value = "Employment (% of total)"
$ws.Range("A12").Font.Bold = $true

$ws.Range("D12").Value = "'50"

# --- Enterprises row (carried down from the old row 6) -------------------
$ws.Range("A13").Value = "Enterprises (% of total)"
$ws.Range("A13").Font.Bold = $true

$ws.Range("D13").Value = "'85"

# --- Source line (carried down from the old row 7) -----------------------
$ws.Range("A14").Value = "Source: MFA, 2010"
$ws.Range("A14").Font.Italic = $true

# --- New "MFA" detail block at the bottom of the sheet --------------------
$ws.Range("A20").Value = "MFA"
$ws.Range("A20").Font.Bold = $true

$ws.Range("A21").Value = 'Ministry of Foreign Affaris (MFA), "SMALL TO MEDIUM ENTERPRISE PAPERS", N/S, p. 3. Available at http://mfa.gov.af/content/files/SME%20PAPER.pdf'
$ws.Range("A21").Font.Italic = $true
